$d = $word.ActiveDocument

# 1) Capitalise "bessons" -> "Bessons" in "Constel·lació de bessons" (5 occurrences
#    throughout the document, including the title/heading and the body paragraph).
$null = $d.Content.Find.Execute(
    "Constel" + [char]0x00B7 + "lació de bessons",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Constel" + [char]0x00B7 + "lació de Bessons", 2)

# 2) Update the astromap year (2019 -> 2022) in the "Jenik Hollan, CzechGlobe (...)."
#    paragraph, and simplify it down to a single plain run (no hyperlink styling),
#    preceded by an empty run - matching the canonical edit.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Jenik Hollan*GaNight*") {
        $target = $para.Range
    }
}

if ($target -ne $null) {
    $rangeEnd = $target.End - 1
    $rr = $d.Range($target.Start, $rangeEnd)
    $newText = "Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
    $xmlFrag = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $null = $rr.InsertXML($xmlFrag)
}
